# Add a new "Login" worksheet ahead of the existing "PythonArray" sheet and
# populate it with a small username/password table (data-driven login test
# data), matching the author's commit "changed TestData for Datadriven Login".

$wb = $excel.ActiveWorkbook

# Remember the pre-existing sheet before we start inserting new ones.
$pythonArray = $wb.Worksheets.Item("PythonArray")

# Worksheets.Add() inserts the new sheet immediately before the ActiveSheet,
# which is still "PythonArray" at this point, so "Login" lands first.
$login = $wb.Worksheets.Add()
$login.Name = "Login"

# Make sure tab order is exactly [Login, PythonArray].
$login.Move($pythonArray, $null)

# Header row.
$login.Range("A1").Value = "Username"
$login.Range("B1").Value = "Password"

# Credential row.
$login.Range("A2").Value = "Group1@dslgo"
$login.Range("B2").Value = "dsalgo2024"

# Style the credential row like a code/monospace cell with a highlighted
# background, matching the rest of the workbook's "code" look.
$cred = $login.Range("A2:B2")
$cred.Font.Name = "Menlo"
$cred.Font.Size = 12
$cred.Font.Color = 16711722
$cred.Interior.Color = 16777215

Write-Host "Added Login sheet with credentials ahead of PythonArray"
